$d = $word.ActiveDocument

# Update the date line (unique text, safe to use Find & Replace)
$d.Content.Find.Execute("2023-11-20 Monday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-11-21 Tuesday", 2)

# Update the division problems in the table. Several values repeat (e.g. the
# text "97÷9=10, 7" occurs twice and must become two different results;
# "99÷6=16, 3" is both an old value and a later new value), so a single
# document-wide Find/Replace is unsafe. Instead address each answer cell
# directly by its (row, column) position in the table.
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text  = "98÷9=10, 8"
$t.Cell(1,2).Range.Text  = "71÷8=8, 7"
$t.Cell(1,3).Range.Text  = "15÷8=1, 7"
$t.Cell(1,4).Range.Text  = "79÷7=11, 2"
$t.Cell(1,5).Range.Text  = "28÷7=4, 0"

$t.Cell(5,1).Range.Text  = "78÷8=9, 6"
$t.Cell(5,2).Range.Text  = "30÷2=15, 0"
$t.Cell(5,3).Range.Text  = "72÷7=10, 2"
$t.Cell(5,4).Range.Text  = "87÷7=12, 3"
$t.Cell(5,5).Range.Text  = "18÷4=4, 2"

$t.Cell(9,1).Range.Text  = "33÷4=8, 1"
$t.Cell(9,2).Range.Text  = "16÷4=4, 0"
$t.Cell(9,3).Range.Text  = "94÷7=13, 3"
$t.Cell(9,4).Range.Text  = "79÷4=19, 3"
$t.Cell(9,5).Range.Text  = "86÷6=14, 2"

$t.Cell(13,1).Range.Text = "78÷2=39, 0"
$t.Cell(13,2).Range.Text = "99÷6=16, 3"
$t.Cell(13,3).Range.Text = "62÷8=7, 6"
$t.Cell(13,4).Range.Text = "17÷9=1, 8"
$t.Cell(13,5).Range.Text = "74÷8=9, 2"

$t.Cell(17,1).Range.Text = "53÷2=26, 1"
$t.Cell(17,2).Range.Text = "20÷6=3, 2"
$t.Cell(17,3).Range.Text = "99÷3=33, 0"
$t.Cell(17,4).Range.Text = "86÷4=21, 2"
$t.Cell(17,5).Range.Text = "20÷5=4, 0"
